# Update column G ("K") values on Sheet1 per the regenerated save_data.
# The commit regenerates the K column (previously "Strike#") with newly
# calculated values; only the literal numbers change, rows/columns stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 3
    7  = 0
    8  = 1
    9  = 1
    10 = 2
    11 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 3
    25 = 1
    26 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
